# Implement loading/saving of error values.
# Adds two new example rows to the "Data Types" sheet:
#   row 21: a literal error value (#VALUE!)
#   row 22: a formula that evaluates to an error (1/0 -> #DIV/0!)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 - error value entered as a literal (not via a formula).
$ws.Range("B21").Value = "Error from literal:"
$ws.Range("C21").Value = "#VALUE!"

# Row 22 - error value produced by evaluating a formula.
$ws.Range("B22").Value = "Error from evaluation:"
$ws.Range("C22").Formula = "=1/0"
